$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays as Text so numeric-looking values are not
# auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.169.04'
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").Value = '1.856.51'
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '232.09'
$ws.Range("E5").Value = '  -3.06%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '0.4727'
$ws.Range("E7").Value = '  -2.04%  '
$ws.Range("D8").Value = '0.2734'
$ws.Range("E8").Value = '  -3.52%  '
$ws.Range("D9").Value = '0.06406'
$ws.Range("E9").Value = '  -1.87%  '
$ws.Range("D10").Value = '1.861.66'
$ws.Range("E10").Value = '  -1.29%  '
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = '16.24'
$ws.Range("E12").Value = '  -2.28%  '
$ws.Range("D13").Value = '5.022'
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").Value = '85.14'
$ws.Range("E14").Value = '  -4.14%  '
$ws.Range("D15").Value = '0.6300'
$ws.Range("E15").Value = '  -5.37%  '
$ws.Range("D16").Value = '30.117.33'
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '230.79'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '12.72'
$ws.Range("E19").Value = '  -4.75%  '
$ws.Range("D20").Value = '0.000007321'
$ws.Range("E20").Value = '  -3.82%  '
$ws.Range("D21").Value = '2.096.89'
$ws.Range("E21").Value = '  -5.20%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Value = '5.049'
$ws.Range("E23").Value = '  -4.64%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '5.988'
$ws.Range("E24").Value = '  -3.11%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.261'
$ws.Range("E25").Value = '  -1.29%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '165.38'
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '17.79'
$ws.Range("E27").Value = '  -4.55%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '1.880'
$ws.Range("E28").Value = '  -3.42%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '0.1011'
$ws.Range("E29").Value = '  +5.80%  '
$ws.Range("E30").Value = '  -2.73%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '4.137'
$ws.Range("E31").Value = '  -5.55%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '3.922'
$ws.Range("E32").Value = '  -2.98%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.04893'
$ws.Range("E33").Value = '  -2.76%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.139'
$ws.Range("E34").Value = '  -5.82%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7204'
$ws.Range("E35").Value = '  -3.74%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '0.9995'
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.692'
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01879'
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.635'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '0.9005'
$ws.Range("E40").Value = '  -1.83%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '1.954'
$ws.Range("E41").Value = '  -6.86%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '105.49'
$ws.Range("E42").Value = '  -0.81%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '0.9991'
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.4089'
$ws.Range("E44").Value = '  -4.72%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '5.538'
$ws.Range("E45").Value = '  -4.44%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.062'
$ws.Range("E46").Value = '  -4.85%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '61.11'
$ws.Range("E47").Value = '  -4.98%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.1197'
$ws.Range("E48").Value = '  -6.98%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.757'
$ws.Range("E49").Value = '  -2.11%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '33.20'
$ws.Range("E50").Value = '  -2.23%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.396'
$ws.Range("E51").Value = '  -5.69%  '
